# Insert two new data rows (new rows 783 and 784) into the weekly Uva
# (grape) price table on Sheet1. All existing rows from the old row 783
# onward shift down by two rows (old 783 -> new 785, ..., old 834 -> new 836).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 783, pushing everything from the old row 783
# down to row 785 (matches xlShiftDown semantics of EntireRow.Insert).
$ws.Range("A783:T784").EntireRow.Insert()

# ---- New row 783: Thompson seedless / Región Metropolitana ----
$ws.Cells.Item(783, 1).Value = 9
$ws.Cells.Item(783, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(783, 3).Value = "Metropolitana"
$ws.Cells.Item(783, 4).Value = 45021
$ws.Cells.Item(783, 5).Value = 13
$ws.Cells.Item(783, 6).Value = "Fruta"
$ws.Cells.Item(783, 7).Value = 100109
$ws.Cells.Item(783, 8).Value = "Uva"
$ws.Cells.Item(783, 9).Value = 100109001
$ws.Cells.Item(783, 10).Value = "Uva"
$ws.Cells.Item(783, 11).Value = "Thompson seedless"
$ws.Cells.Item(783, 12).Value = "Primera"
$ws.Cells.Item(783, 13).Value = 350
$ws.Cells.Item(783, 14).Value = 12000
$ws.Cells.Item(783, 15).Value = 12000
$ws.Cells.Item(783, 16).Value = 12000
$ws.Cells.Item(783, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(783, 18).Value = "Región Metropolitana"
$ws.Cells.Item(783, 19).Value = 667
$ws.Cells.Item(783, 20).Value = 18

# ---- New row 784: Timco / Región de O'Higgins ----
$ws.Cells.Item(784, 1).Value = 9
$ws.Cells.Item(784, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(784, 3).Value = "Metropolitana"
$ws.Cells.Item(784, 4).Value = 45021
$ws.Cells.Item(784, 5).Value = 13
$ws.Cells.Item(784, 6).Value = "Fruta"
$ws.Cells.Item(784, 7).Value = 100109
$ws.Cells.Item(784, 8).Value = "Uva"
$ws.Cells.Item(784, 9).Value = 100109001
$ws.Cells.Item(784, 10).Value = "Uva"
$ws.Cells.Item(784, 11).Value = "Timco"
$ws.Cells.Item(784, 12).Value = "Primera"
$ws.Cells.Item(784, 13).Value = 300
$ws.Cells.Item(784, 14).Value = 12000
$ws.Cells.Item(784, 15).Value = 12000
$ws.Cells.Item(784, 16).Value = 12000
$ws.Cells.Item(784, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(784, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(784, 19).Value = 667
$ws.Cells.Item(784, 20).Value = 18
